$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "About": update source attribution from NREL RE Futures to
# the 2035 Report (GSPP, GridLab, Energy Innovation), and drop the old
# "80% RE-ITI (2014)" scenario note.
# ---------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

$about.Range("B13").Value = "GSPP, GridLab, Energy Innovation"
$about.Range("B14").Value = 2020
$about.Range("B15").Value = "The 2035 Report"
$about.Range("B16").Value = "https://www.2035report.com/data-explorer/"

# Remove the hyperlink that used to point at the NREL data viewer (B16)
# while leaving the other two hyperlinks (B9, B23) untouched.
foreach ($h in $about.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$B$16') {
        $h.Delete()
    }
}

# The old "80% RE-ITI (2014)" scenario footnote no longer applies.
$about.Range("B17").ClearContents()

# ---------------------------------------------------------------------
# Sheet "RE Futures Data & Calcs": swap the NREL RE-Futures 80% RE-ITI
# numbers for the 2035 Report figures.
# ---------------------------------------------------------------------
$calcs = $wb.Worksheets.Item("RE Futures Data & Calcs")

# End year for the new scenario.
$calcs.Range("B2").Value = 2035

# Start-year battery capacity now pulled live from BGBSC!R2 (GW), rather
# than the old hard-coded 26.
$calcs.Range("B3").Formula = "=BGBSC!R2/10^3"
$calcs.Range("B3").NumberFormat = "0"

# Scenario label + end-year target capacity (GW).
$calcs.Range("A4").Value = "2035 Scenario"
$calcs.Range("B4").Value = 150

# End-of-trend capacity target is now a hard 500 GW figure rather than a
# reference back to B4.
$calcs.Range("C15").Value = 500

# Relabel the annualized-trend row.
$calcs.Range("A18").Value = "2034 Report Annualized"

# ---------------------------------------------------------------------
# Make PAGBSC the active sheet/tab, matching the saved workbook state.
# ---------------------------------------------------------------------
$pagbsc = $wb.Worksheets.Item("PAGBSC")
$pagbsc.Activate()
